$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.700207233428955
$ws.Range("B1").Value = 3.863176107406616
$ws.Range("C1").Value = 4.462451934814453
$ws.Range("D1").Value = 2.431181669235229
$ws.Range("E1").Value = 1.575451374053955
